$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row to the "Form_Responses" table (grows the table range from A1:Q4 to A1:Q5)
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Seed the new row (row 5) with the same cell formatting as row 3, which already
# carries the exact style pattern needed for this row (date, text, quote-prefixed
# "text that looks like a number" for phone fields, and hyperlink-styled photo columns)
$ws.Range("A3:Q3").Copy()
$ws.Range("A5:Q5").PasteSpecial(-4122)

# Populate the new response row - "Jóvenes de Acción Católica Mercedaria"
$ws.Cells.Item(5, 1).Value = 46013.561763125
$ws.Cells.Item(5, 2).Value = "mandreaef@gmail.com"
$ws.Cells.Item(5, 3).Value = "Jóvenes de Acción Católica Mercedaria"
$ws.Cells.Item(5, 4).Value = "@jacm.sanramon"
$ws.Cells.Item(5, 5).Value = "San Ramón Nonato"
$ws.Cells.Item(5, 6).Value = "Pastoral Juvenil Parroquial adherida a los estatutos de la Acción Católica (asociación de laicos para laicos) y siguiendo el carisma Mercedario (carisma de la Parroquia)"
$ws.Cells.Item(5, 7).Value = "- Encuentros Formativos: Espacios de formación para los jóvenes de la parroquia (pertenecientes al grupo e invitados)`n- Liturgia: Encargados de la liturgia de la Eucaristía de los Viernes a las 5:00 PM (previa a los encuentros)`n- Acción Social: Por los menos dos o tres acciones sociales trimestrales en los barrios pertenecientes a la comunidad`n- Experiencias Espirituales: Convivencias y retiros espirituales abiertos para jóvenes de la parroquia`n- JAC-MEM: Encuentros combinados con Pastoral Familiar como método de preparación para la generación de relevo"
$ws.Cells.Item(5, 8).Value = "Victor Kneider"
$ws.Cells.Item(5, 9).Value = "'04246501227"
$ws.Cells.Item(5, 10).Value = "Maria Andrea Espina"
$ws.Cells.Item(5, 11).Value = "'04146019600"
$ws.Cells.Item(5, 12).Value = "Jóvenes"
$ws.Cells.Item(5, 13).Value = "No"
$ws.Cells.Item(5, 14).Value = "Viernes"
$ws.Cells.Item(5, 15).Value = "Viernes de 5 PM a 8 PM"
$ws.Cells.Item(5, 16).Value = "https://drive.google.com/open?id=1Tw8D3hctSFlOJIcZ7JQLBNhUrNlGPpea"
$ws.Cells.Item(5, 17).Value = "https://drive.google.com/open?id=1LLWyeQkEfAM1vR32udBDK4_ija6gUjSg"

# Turn the two photo-link cells into real hyperlinks (adds rows to <hyperlinks> and
# the accompanying external relationships), same as the other response rows
$ws.Hyperlinks.Add($ws.Range("P5"), "https://drive.google.com/open?id=1Tw8D3hctSFlOJIcZ7JQLBNhUrNlGPpea")
$ws.Hyperlinks.Add($ws.Range("Q5"), "https://drive.google.com/open?id=1LLWyeQkEfAM1vR32udBDK4_ija6gUjSg")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" style; restore
# the existing custom hyperlink-look formatting used by the rest of the table instead
$ws.Range("P3:Q3").Copy()
$ws.Range("P5:Q5").PasteSpecial(-4122)
$wb.Styles.Item("Hyperlink").Delete()

# Match the row height used by the rest of the table
$ws.Rows.Item(5).RowHeight = 22.5
